$wb = $excel.ActiveWorkbook

# --- Sheet: Significant Components ---
$ws1 = $wb.Worksheets.Item("Significant Components")
$ws1.Range("C2").Value = "['QEXTRCT' 'QEDLESHI' 'QESL' 'QNOHLTH' 'PPUNIT' 'QSERV' 'QHISPC' 'QFHH'`n 'PERCAP']"
$ws1.Range("C4").Value = "['MEDAGE' 'QSSBEN' 'QAGEDEP']"
$ws1.Range("C6").Value = "['PPUNIT' 'QRENTER' 'QNOAUTO' 'QPOVTY']"

# --- Sheet: Loading Factors ---
$ws2 = $wb.Worksheets.Item("Loading Factors")
$lf = New-Object 'object[,]' 19,6
$lf[0,0] = "QEXTRCT"
$lf[0,1] = 0.7677512013307708
$lf[0,2] = 0.1449653481766418
$lf[0,3] = 0.01129679077873671
$lf[0,4] = -0.2382017695991725
$lf[0,5] = 0.09278788174314587
$lf[1,0] = "QEDLESHI"
$lf[1,1] = 0.8777939278338194
$lf[1,2] = 0.213077928978771
$lf[1,3] = -0.01839361500093432
$lf[1,4] = -0.1076883429731288
$lf[1,5] = 0.1846338367566293
$lf[2,0] = "QESL"
$lf[2,1] = 0.8009661475849839
$lf[2,2] = 0.1517950795169672
$lf[2,3] = -0.03424774416423925
$lf[2,4] = -0.2374898287916476
$lf[2,5] = 0.2037776713559489
$lf[3,0] = "QNOHLTH"
$lf[3,1] = 0.6889886573070284
$lf[3,2] = 0.414959225541808
$lf[3,3] = -0.1190897065630508
$lf[3,4] = -0.1154839642310589
$lf[3,5] = 0.2786568866560452
$lf[4,0] = "PPUNIT"
$lf[4,1] = 0.7302277711183555
$lf[4,2] = -0.004620272024895255
$lf[4,3] = -0.1512975164200612
$lf[4,4] = 0.05458730556150756
$lf[4,5] = -0.4752714593490717
$lf[5,0] = "QSERV"
$lf[5,1] = 0.581704943683128
$lf[5,2] = 0.3577287418172193
$lf[5,3] = -0.2244058884840431
$lf[5,4] = -0.03201466441996235
$lf[5,5] = 0.2754259447697267
$lf[6,0] = "QHISPC"
$lf[6,1] = 0.8328587063035261
$lf[6,2] = 0.3339037468887154
$lf[6,3] = -0.1364957918037546
$lf[6,4] = -0.126970430831828
$lf[6,5] = 0.09670564462395792
$lf[7,0] = "QFHH"
$lf[7,1] = 0.5630560183163126
$lf[7,2] = 0.3008187155964411
$lf[7,3] = -0.09551317755867508
$lf[7,4] = 0.2634165381813704
$lf[7,5] = -0.03178270184700521
$lf[8,0] = "PERCAP"
$lf[8,1] = 0.4895374235147206
$lf[8,2] = 0.7214540992653554
$lf[8,3] = -0.2685725828240927
$lf[8,4] = 0.05482090093234365
$lf[8,5] = 0.1831253412505047
$lf[9,0] = "QRICH"
$lf[9,1] = 0.2150971482753866
$lf[9,2] = 0.870130548924468
$lf[9,3] = -0.1729034375056023
$lf[9,4] = -0.0142515353605643
$lf[9,5] = 0.2948729845659209
$lf[10,0] = "MDHSEVAL"
$lf[10,1] = 0.3857664718237477
$lf[10,2] = 0.8013370028956395
$lf[10,3] = -0.03601197413356935
$lf[10,4] = -0.02873185735955726
$lf[10,5] = -0.03016150937756498
$lf[11,0] = "MEDAGE"
$lf[11,1] = -0.3105472617371984
$lf[11,2] = -0.2465384280778078
$lf[11,3] = 0.7910361908750342
$lf[11,4] = -0.01289474291113122
$lf[11,5] = -0.2717085724916322
$lf[12,0] = "QSSBEN"
$lf[12,1] = 0.01836498286088173
$lf[12,2] = -0.05367045675897775
$lf[12,3] = 0.7773306962764825
$lf[12,4] = 0.136209976584241
$lf[12,5] = -0.1455945068606067
$lf[13,0] = "QAGEDEP"
$lf[13,1] = -0.04282432190847357
$lf[13,2] = -0.1184872960740083
$lf[13,3] = 0.6543380961981018
$lf[13,4] = 0.6427560452936805
$lf[13,5] = -0.1139548183463203
$lf[14,0] = "QFEMLBR"
$lf[14,1] = -0.2416846133880868
$lf[14,2] = 0.08178103263001013
$lf[14,3] = -0.02959553031661272
$lf[14,4] = 0.7849929020395576
$lf[14,5] = 0.003645848475871107
$lf[15,0] = "QFEMALE"
$lf[15,1] = -0.04556236170050956
$lf[15,2] = -0.05659266582526922
$lf[15,3] = 0.1671554332667903
$lf[15,4] = 0.8778431990123544
$lf[15,5] = -0.0242050142911791
$lf[16,0] = "QRENTER"
$lf[16,1] = 0.01759024404216022
$lf[16,2] = 0.2288188101098284
$lf[16,3] = -0.4234144374084097
$lf[16,4] = -0.09662324389171931
$lf[16,5] = 0.7659446219419647
$lf[17,0] = "QNOAUTO"
$lf[17,1] = 0.1660837693831291
$lf[17,2] = 0.06295541374020226
$lf[17,3] = -0.1064187480161447
$lf[17,4] = -0.01530325772402283
$lf[17,5] = 0.6312035357109876
$lf[18,0] = "QPOVTY"
$lf[18,1] = 0.3701334443131513
$lf[18,2] = 0.157633475077385
$lf[18,3] = -0.3817953811110359
$lf[18,4] = 0.08010853240336936
$lf[18,5] = 0.4611334541988245
$ws2.Range("A2:F20").Value = $lf

# --- Sheet: All Refactor Variances ---
$ws3 = $wb.Worksheets.Item("All Refactor Variances")
$arv = New-Object 'object[,]' 4,10
$arv[0,0] = 4.839746425340321
$arv[0,1] = 3.422021238437344
$arv[0,2] = 2.232417334146756
$arv[0,3] = 2.056337886114042
$arv[0,4] = 2.044708471301834
$arv[0,5] = 5.117234792812528
$arv[0,6] = 2.680150619567169
$arv[0,7] = 2.243476605342225
$arv[0,8] = 2.066288181755519
$arv[0,9] = 1.901714776046654
$arv[1,0] = 0.2304641154923962
$arv[1,1] = 0.1629533923065402
$arv[1,2] = 0.1063055873403217
$arv[1,3] = 0.09792085171971628
$arv[1,4] = 0.09736707006199209
$arv[1,5] = 0.2693281469901331
$arv[1,6] = 0.1410605589245879
$arv[1,7] = 0.1180777160706434
$arv[1,8] = 0.10875200956608
$arv[1,9] = 0.1000902513708765
$arv[2,0] = 0.2304641154923962
$arv[2,1] = 0.3934175077989364
$arv[2,2] = 0.4997230951392581
$arv[2,3] = 0.5976439468589744
$arv[2,4] = 0.6950110169209665
$arv[2,5] = 0.2693281469901331
$arv[2,6] = 0.4103887059147209
$arv[2,7] = 0.5284664219853643
$arv[2,8] = 0.6372184315514443
$arv[2,9] = 0.7373086829223208
$arv[3,0] = 0.3315977874903292
$arv[3,1] = 0.2344615960599521
$arv[3,2] = 0.1529552550278643
$arv[3,3] = 0.1408910784659567
$arv[3,4] = 0.1400942829558977
$arv[3,5] = 0.3652854675773677
$arv[3,6] = 0.1913181848957682
$arv[3,7] = 0.1601469219142283
$arv[3,8] = 0.1474986150102582
$arv[3,9] = 0.1357508106023777
$ws3.Range("I2:R5").Value = $arv

# --- Sheet: Final Variances ---
$ws4 = $wb.Worksheets.Item("Final Variances")
$fv = New-Object 'object[,]' 4,5
$fv[0,0] = 5.117234792812528
$fv[0,1] = 2.680150619567169
$fv[0,2] = 2.243476605342225
$fv[0,3] = 2.066288181755519
$fv[0,4] = 1.901714776046654
$fv[1,0] = 0.2693281469901331
$fv[1,1] = 0.1410605589245879
$fv[1,2] = 0.1180777160706434
$fv[1,3] = 0.10875200956608
$fv[1,4] = 0.1000902513708765
$fv[2,0] = 0.2693281469901331
$fv[2,1] = 0.4103887059147209
$fv[2,2] = 0.5284664219853643
$fv[2,3] = 0.6372184315514443
$fv[2,4] = 0.7373086829223208
$fv[3,0] = 0.3652854675773677
$fv[3,1] = 0.1913181848957682
$fv[3,2] = 0.1601469219142283
$fv[3,3] = 0.1474986150102582
$fv[3,4] = 0.1357508106023777
$ws4.Range("B2:F5").Value = $fv

# --- Sheet: Included and Excluded ---
$ws5 = $wb.Worksheets.Item("Included and Excluded")
$ws5.Range("B2").Value = "[['QEXTRCT', 'QEDLESHI', 'QESL', 'QNOHLTH', 'PPUNIT', 'QSERV', 'QHISPC', 'QFHH', 'PERCAP', 'QRICH', 'MDHSEVAL', 'MEDAGE', 'QSSBEN', 'QAGEDEP', 'QFEMLBR', 'QFEMALE', 'QRENTER', 'QNOAUTO', 'QPOVTY']]"
